$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40 - pushes existing rows 40..91 down to 41..92
$ws.Rows("40:40").Insert()

# Populate the newly inserted row 40 with the new weekly data point
$ws.Range("A40").Value = 11
$ws.Range("B40").Value = "Vega Monumental Concepción"
$ws.Range("C40").Value = "Bíobío"
$ws.Range("D40").Value = 44580
$ws.Range("E40").Value = 8
$ws.Range("F40").Value = 100112032
$ws.Range("G40").Value = "Zapallo italiano"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 100
$ws.Range("K40").Value = 14000
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = 14500
$ws.Range("N40").Value = "$/caja 50 unidades"
$ws.Range("O40").Value = "Región de O'Higgins"
$ws.Range("P40").Value = 290
$ws.Range("Q40").Value = 50
$ws.Range("R40").Value = "Hortaliza"
